# Locate the paragraph that contains the "pip install name_module ... установка модуля"
# line (center-justified, bold command + translation), so we can clone it into a new
# paragraph right below with "uninstall"/"удаление" wording, as shown in the diff.
$d = $word.ActiveDocument

$sourceIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $paraText = $d.Paragraphs.Item($i).Range.Text
    if ($paraText -like "*pip install name_module*" -and $paraText -like "*установка*") {
        $sourceIndex = $i
    }
}

$sourcePara = $d.Paragraphs.Item($sourceIndex)

# Duplicate the whole paragraph (including its own paragraph mark, so the new copy
# gets its own <w:p> with the same pPr) right after the source paragraph.
$sourceRange = $d.Range($sourcePara.Range.Start, $sourcePara.Range.End)
$clonedFormattedText = $sourceRange.FormattedText
$insertionPoint = $d.Range($sourcePara.Range.End, $sourcePara.Range.End)
$insertionPoint.FormattedText = $clonedFormattedText

$newPara = $d.Paragraphs.Item($sourceIndex + 1)
$newParaStart = $newPara.Range.Start

# "pip install " and "name_module" share identical run formatting, so a plain
# Find/Replace on "install" would merge them into a single run. Drop a throwaway
# bookmark at that run boundary first -- it keeps the two runs apart through the
# edit, then gets deleted, leaving two separate <w:r> elements (as in the target
# XML) with no bookmark residue.
$boundary = $newParaStart + [int]"pip install ".Length
$d.Bookmarks.Add("zzTmpRunSplit", $d.Range($boundary, $boundary)) | Out-Null

$newPara.Range.Find.Execute("install", $true, $false, $false, $false, $false, $true, 1, $false, "uninstall", 2) | Out-Null
$d.Paragraphs.Item($sourceIndex + 1).Range.Find.Execute("установка", $true, $false, $false, $false, $false, $true, 1, $false, "удаление", 2) | Out-Null

$d.Bookmarks.Item("zzTmpRunSplit").Delete()
